$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: bring back weather data values
$ws.Range("F10").Value = 15649
$ws.Range("K10").Value = 41595
$ws.Range("O10").Value = 15650.5649
$ws.Range("T10").Value = 41599.1595
